$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.099.42"
$ws.Range("E2").Value = "  -4.71%  "

$ws.Range("D3").Value = "3.288.22"
$ws.Range("E3").Value = "  -5.44%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.88"
$ws.Range("E5").Value = "  -3.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.92"
$ws.Range("E6").Value = "  -3.80%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.286.66"
$ws.Range("E8").Value = "  -5.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("E9").Value = "  -1.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.28"
$ws.Range("E10").Value = "  -4.66%  "

$ws.Range("E11").Value = "  -4.61%  "

$ws.Range("E12").Value = "  -3.48%  "

$ws.Range("D13").Value = "3.856.70"
$ws.Range("E13").Value = "  -5.35%  "

$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").Value = "3.298.31"
$ws.Range("E15").Value = "  -5.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000166"
$ws.Range("E16").Value = "  -6.28%  "

$ws.Range("D17").Value = "61.191.54"
$ws.Range("E17").Value = "  -4.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "23.99"
$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.57"
$ws.Range("E19").Value = "  -2.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.18"
$ws.Range("E20").Value = "  -2.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.87"
$ws.Range("E21").Value = "  -10.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "351.19"
$ws.Range("E22").Value = "  -8.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.549"
$ws.Range("E23").Value = "  -4.42%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").Value = "3.424.96"
$ws.Range("E25").Value = "  -5.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "68.88"
$ws.Range("E26").Value = "  -7.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000106"
$ws.Range("E27").Value = "  -5.56%  "

$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.07"
$ws.Range("E29").Value = "  -1.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.43"
$ws.Range("E30").Value = "  -0.53%  "

$ws.Range("E31").Value = "  -6.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.73"
$ws.Range("E32").Value = "  -2.71%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("E34").Value = "  -3.35%  "

$ws.Range("D35").Value = "3.320.62"
$ws.Range("E35").Value = "  -5.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.47"
$ws.Range("E36").Value = "  -2.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.20"
$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.73"
$ws.Range("E38").Value = "  -0.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "163.15"
$ws.Range("E39").Value = "  +0.18%  "

$ws.Range("E40").Value = "  -3.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0747"
$ws.Range("E41").Value = "  -3.55%  "

$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.10"
$ws.Range("E43").Value = "  -0.86%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.36"
$ws.Range("E44").Value = "  +1.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.737"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.11"
$ws.Range("E46").Value = "  -2.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.54"
$ws.Range("E47").Value = "  -5.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.03"
$ws.Range("E48").Value = "  -7.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.63"
$ws.Range("E49").Value = "  -1.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.845"
$ws.Range("E50").Value = "  -7.71%  "

$ws.Range("E51").Value = "  +2.47%  "
